$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-by-cell updates per the source diff. Numeric-looking price strings
# (single decimal point, e.g. "205.99") are prefixed with a leading
# apostrophe so Excel stores them as text, matching the original
# inlineStr/text representation instead of silently converting them to
# floating point numbers.

$ws.Range('D2').Value = '26.886.62'
$ws.Range('D3').Value = '1.565.24'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''205.99'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '''21.77'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').Value = '1.788.20'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '1.580.31'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '''3.73'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '26.894.04'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').Value = '''61.34'
$ws.Range('E17').Value = '  -2.72%  '
$ws.Range('D18').Value = '''215.27'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').Value = '''7.38'
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').Value = '''9.21'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').Value = '''1.99'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Value = '''154.18'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('E26').Value = '  +1.74%  '
$ws.Range('D27').Value = '''14.95'
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  -3.32%  '
$ws.Range('D32').Value = '''3.16'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').Value = '1.404.51'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('D35').Value = '''1.52'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').Value = '''0.921'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').Value = '''0.998'
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('E43').Value = '  +2.64%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '''2.18'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''1.76'
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('D47').Value = '1.701.03'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').Value = '''86.31'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('D49').Value = '''0.0506'
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').Value = '0.0₇0983'
$ws.Range('E50').Value = '  -1.44%  '
$ws.Range('E51').Value = '  +0.28%  '
